# Update the "想去人数" (F column) counts on the "展览" sheet and on the
# "全部类型" sheet (which mirrors the same events at a one-row offset).
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - rId1 / sheet1.xml
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F5").Value  = 5204
$wsExhibition.Range("F9").Value  = 583
$wsExhibition.Range("F11").Value = 1046
$wsExhibition.Range("F13").Value = 1438
$wsExhibition.Range("F14").Value = 4145
$wsExhibition.Range("F15").Value = 429
$wsExhibition.Range("F17").Value = 147
$wsExhibition.Range("F19").Value = 3179
$wsExhibition.Range("F20").Value = 156
$wsExhibition.Range("F21").Value = 1061
$wsExhibition.Range("F25").Value = 98
$wsExhibition.Range("F26").Value = 27
$wsExhibition.Range("F27").Value = 134
$wsExhibition.Range("F29").Value = 288
$wsExhibition.Range("F30").Value = 22
$wsExhibition.Range("F31").Value = 51
$wsExhibition.Range("F32").Value = 11
$wsExhibition.Range("F34").Value = 11

# Sheet "全部类型" (All types) - rId4 / sheet4.xml
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value  = 5205
$wsAll.Range("F10").Value = 583
$wsAll.Range("F12").Value = 1046
$wsAll.Range("F14").Value = 1438
$wsAll.Range("F15").Value = 4145
$wsAll.Range("F16").Value = 429
$wsAll.Range("F18").Value = 147
$wsAll.Range("F20").Value = 3179
$wsAll.Range("F21").Value = 156
$wsAll.Range("F22").Value = 1061
$wsAll.Range("F26").Value = 98
$wsAll.Range("F27").Value = 27
$wsAll.Range("F28").Value = 134
$wsAll.Range("F30").Value = 288
$wsAll.Range("F31").Value = 22
$wsAll.Range("F32").Value = 51
$wsAll.Range("F33").Value = 11
$wsAll.Range("F35").Value = 11

$wb.Save()
